$d = $word.ActiveDocument

# --- 1. The (single, empty) body paragraph now carries the "Раздел"
#        ("Chapter/Section title") paragraph style, with an explicit
#        "no numbering" override (ilvl 0 / numId 0) and a 709-twip
#        (1.25 cm) left indent. ---
$p = $d.Paragraphs(1)
$p.Style = "a3"
$p.Range.ParagraphFormat.LeftIndent = 35.45   # 35.45 pt == 709 twips
$p.Range.ListFormat.RemoveNumbers()

# --- 2. Style "a3" itself: renamed from "Глава" to "Раздел", and now
#        forces a page break before every paragraph using it, so that
#        chapter/section titles always start on a new page. ---
$s = $d.Styles("a3")
$s.NameLocal = "Раздел"
$s.ParagraphFormat.PageBreakBefore = $true
